# Apply the forensic-case worksheet corrections:
#  - fix a typo in the Ex1 notes ("Al Copne" -> "Al Copone")
#  - fill in the previously-blank summary/reasonForRemoval cells for exhibit row 3 (Ex9)
#  - clear the (no longer applicable) imagingType for exhibit row 3 (Ex9)
#  - flip exportedEvidence for exhibit row 3 (Ex9) from Y to N
#  - fill in the previously-blank imagingType for exhibit row 5 (Ex12)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the notes for exhibit row 2 (Ex1): "Al Copne" -> "Al Copone"
$ws.Range("I2").Value = "The operating system was Android 11.This was from Al Copone's phone and he gave 193104 as the pincode. This phone contained 176 instant messages between AlCopone@gmail.com and 16185551211 (James). There were message between Al Copone and his accountant."

# Exhibit row 3 (Ex9 - DVR): fill in summary, reasonForRemoval; clear imagingType; exportedEvidence Y -> N
$ws.Range("J3").Value = "On July 28, 1931, Sherlock Holmes attended the warrant at the 7244 Prairie Avenue, Chicago Illinois."
$ws.Range("J3").HorizontalAlignment = 1
$ws.Range("AJ3").Value = "finished"
$ws.Range("AO3").Value = ""
$ws.Range("BD3").Value = "N"

# Exhibit row 5 (Ex12 - vehicle): fill in imagingType
$ws.Range("AO5").Value = "advanced logical"

# Leave the cursor where the editor left it
$null = $ws.Range("H10").Select()
